$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '71.551.12'
$ws.Range('E2').Value = '  +3.26%  '
$ws.Range('D3').Value = '3.706.64'
$ws.Range('E3').Value = '  +8.48%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '583.18'
$ws.Range('E5').Value = '  +0.64%  '
Set-TextValue 'D6' '179.10'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('D7').Value = '3.691.07'
$ws.Range('E7').Value = '  +8.22%  '
Set-TextValue 'D8' '0.616'
$ws.Range('E8').Value = '  +4.37%  '
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('E10').Value = '  +2.03%  '
Set-TextValue 'D12' '49.32'
$ws.Range('E12').Value = '  +1.23%  '
Set-TextValue 'D13' '0.0000287'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').Value = '4.305.41'
$ws.Range('E14').Value = '  +8.52%  '
Set-TextValue 'D15' '683.66'
$ws.Range('E15').Value = '  -1.52%  '
Set-TextValue 'D16' '9.02'
$ws.Range('E16').Value = '  +4.65%  '
$ws.Range('D17').Value = '3.719.01'
$ws.Range('E17').Value = '  +8.81%  '
$ws.Range('D18').Value = '71.690.92'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('E21').Value = '  +2.27%  '
$ws.Range('E22').Value = '  +18.61%  '
Set-TextValue 'D23' '0.944'
Set-TextValue 'D24' '17.50'
$ws.Range('E24').Value = '  +3.66%  '
Set-TextValue 'D25' '102.48'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('E26').Value = '  +3.24%  '
$ws.Range('E27').Value = '  +6.83%  '
Set-TextValue 'D28' '10.41'
$ws.Range('E28').Value = '  +8.83%  '
Set-TextValue 'D29' '35.36'
$ws.Range('E29').Value = '  +6.24%  '
Set-TextValue 'D30' '9.19'
$ws.Range('E30').Value = '  +5.42%  '
Set-TextValue 'D31' '7.33'
$ws.Range('E31').Value = '  +5.96%  '
$ws.Range('E32').Value = '  +10.98%  '
Set-TextValue 'D33' '592.68'
$ws.Range('E33').Value = '  +4.19%  '
Set-TextValue 'D34' '11.21'
$ws.Range('E34').Value = '  +2.14%  '
Set-TextValue 'D35' '0.108'
$ws.Range('E35').Value = '  +4.54%  '
Set-TextValue 'D36' '59.10'
$ws.Range('E36').Value = '  +1.63%  '
Set-TextValue 'D37' '1.00'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').Value = '3.678.05'
$ws.Range('E38').Value = '  +2.73%  '
$ws.Range('E39').Value = '  +5.18%  '
$ws.Range('D40').Value = '0.0₃0770'
$ws.Range('E40').Value = '  +6.41%  '
Set-TextValue 'D41' '35.50'
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('E42').Value = '  +5.77%  '
Set-TextValue 'D43' '2.80'
$ws.Range('E43').Value = '  +5.74%  '
$ws.Range('E44').Value = '  +10.33%  '
$ws.Range('E45').Value = '  +5.51%  '
Set-TextValue 'D46' '2.89'
$ws.Range('E46').Value = '  +9.56%  '
Set-TextValue 'D47' '3.39'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +4.10%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -0.03%  '
Set-TextValue 'D51' '135.79'
$ws.Range('E51').Value = '  +3.01%  '
